$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

$data = @(
  @("Ana Oliveira", 160, 35, 10000),
  @("Bruno Lima", 140, 40, 9500),
  @("Carla Martins", 155, 38, 11200),
  @("Daniel Souza", 170, 36, 12500),
  @("Fernanda Silva", 150, 42, 10800),
  @("Gabriel Santos", 165, 37, 11900),
  @("Helena Rocha", 145, 39, 9800),
  @("Igor Ferreira", 175, 41, 13000),
  @("Julia Almeida", 160, 34, 10200),
  @("Lucas Correia", 150, 36, 9700)
)

$row = 2
foreach ($r in $data) {
  $ws.Cells.Item($row, 1).Value = $r[0]
  $ws.Cells.Item($row, 2).Value = $r[1]
  $ws.Cells.Item($row, 3).Value = $r[2]
  $ws.Cells.Item($row, 4).Value = $r[3]
  $row++
}

$tbl.Resize($ws.Range("A1:D11"))
$tbl.TableStyle = "TableStyleMedium4"

$ws.Range("C2:D11").NumberFormat = """R$"" #,##0.00;[Red]-""R$"" #,##0.00"
